$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($paragraph, [string]$innerP) {
    $xml = $pkgHeader + $innerP + $pkgFooter
    $paragraph.Range.InsertXML($xml)
}

# --- Locate the three target paragraphs by their distinctive text ---
$pNote = $null
$pAndThus = $null
$pValidation = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Note: the opentech server*") { $pNote = $i }
    if ($t -like "And thus the site redirect*") { $pAndThus = $i }
    if ($t -like "You will not be able to see validation output*") { $pValidation = $i }
}

$frag13 = @'
<w:p w14:paraId="669645B0" w14:textId="77777777" w:rsidR="0031044E" w:rsidRDefault="00B57606" w:rsidP="00CD3E3D"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Note: the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>opentech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> server </w:t></w:r><w:r w:rsidR="00673B53"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>has issues</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, likely out of date PHP</w:t></w:r><w:r w:rsidR="0031044E"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p>
'@

$frag17 = @'
<w:p w14:paraId="0DFEB006" w14:textId="2C69074D" w:rsidR="00B57606" w:rsidRPr="0031044E" w:rsidRDefault="0031044E" w:rsidP="00CD3E3D"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:r w:rsidR="00B57606"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">nd </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00B57606"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>thus</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00B57606"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the site redirect do not work correctly despite following all instructions and hours of debugging to remove whitespace and so on.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> The stated workaround in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>the Week 11</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 4_file_redirect slides </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>does not work</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> on the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>opentech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> server.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Neither does removing all whitespace from the code before the header() call or wrapping the call in its own &lt;?</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&gt; tags.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I have also validated via PuTTY that the database on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>opentech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is updated, it is not a database issue.</w:t></w:r></w:p>
'@

$frag19 = @'
<w:p w14:paraId="7428432C" w14:textId="6DEDCBF0" w:rsidR="00B57606" w:rsidRPr="00CD3E3D" w:rsidRDefault="00B57606" w:rsidP="00CB5271"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You will not be able to see validation output due to the issues on the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>opentech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> server, I will prioritize adding them to my demonstration video.</w:t></w:r></w:p>
'@

Set-ParagraphXml $d.Paragraphs($pNote) $frag13
Set-ParagraphXml $d.Paragraphs($pAndThus) $frag17
Set-ParagraphXml $d.Paragraphs($pValidation) $frag19

Write-Output "done"
